$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column (H1), matching the style/formatting of the
# existing header row (e.g. G1 "sum"), and a value of 1 for the first
# data row (H2).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H2").Value = 1
